$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Title slide: "Database Final Project" -> "Movie Database"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Movie Database"

# ---------------------------------------------------------------------------
# 2. Overview slide: split the single paragraph into three paragraphs and
#    append a new closing paragraph.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$overviewBody = $s2.Shapes.Item(2).TextFrame.TextRange
$overviewBody.Text = "Our program acts as a movie compendium that allows clients to search for movies based on many different aspects.'`rIt will allow clients to search for movies based on actors, directors, genre, and more. They can also see reviews on movies and sort movies by rating.`rAfter finding a certain movie, you can then select that movie and see all of the information pertaining to that movie.`r"

# ---------------------------------------------------------------------------
# 3. Duplicate the "Technical Details" slide (index 3) to create the new
#    "Overview Cont." slide, move it right after the Overview slide, and
#    replace its text.
# ---------------------------------------------------------------------------
$dup1 = $p.Slides.Item(3).Duplicate()
$overviewCont = $dup1.Item(1)
$overviewCont.MoveTo(3)

$overviewContTitle = $overviewCont.Shapes.Item(1).TextFrame.TextRange
$overviewContTitle.Text = "Overview Cont."

$overviewContBody = $overviewCont.Shapes.Item(2).TextFrame.TextRange
$overviewContBody.Text = "Users can also add, delete, or update the movies in the database.`rThese changes can affect any aspect of the movie, from title to actors."

$overviewContBodyShape = $overviewCont.Shapes.Item(2)
$overviewContBodyShape.Left = 53.33339582677166
$overviewContBodyShape.Top = 248.84212598425196
$overviewContBodyShape.Width = 676.9029921259843
$overviewContBodyShape.Height = 226.8556792913386

# ---------------------------------------------------------------------------
# 4. "Technical Details" slide (now at index 4) gets its body filled in.
# ---------------------------------------------------------------------------
$techDetails = $p.Slides.Item(4)
$techBody = $techDetails.Shapes.Item(2).TextFrame.TextRange
$techBody.Text = "We wanted to implement the database in a way that was user friendly, but was easy to implement.`rWithout a complex search algorithm like what IMDb has that can search, for example, both actors and movies, we needed to have a way to search for certain things individually.`rWe decided to use a checkbox system that allows the user to select how they are searching for a movie.`r"

# ---------------------------------------------------------------------------
# 5. "Target Users" slide (index 5): extend the final paragraph.
# ---------------------------------------------------------------------------
$targetUsers = $p.Slides.Item(5)
$tuBody = $targetUsers.Shapes.Item(2).TextFrame.TextRange
$lastPara = $tuBody.Paragraphs($tuBody.Paragraphs().Count, 1)
$lastPara.Text = "This program allows users to easily lookup movies through many different ways."

# ---------------------------------------------------------------------------
# 6. Duplicate the "Target Users" slide to create "Future Expansions" at the
#    end of the deck.
# ---------------------------------------------------------------------------
$dup2 = $targetUsers.Duplicate()
$futureExp = $dup2.Item(1)
$futureExp.MoveTo($p.Slides.Count)

$futureExpTitle = $futureExp.Shapes.Item(1).TextFrame.TextRange
$futureExpTitle.Text = "Future Expansions"

$futureExpBody = $futureExp.Shapes.Item(2).TextFrame.TextRange
$futureExpBody.Text = "Allow the ability to search through any combination of search parameters.`rAdd more error handling for incorrect querying.`rAdd the ability to create a favorites list for movies and actors.`rShow more information about the movie, like including a poster of the movie.`rAdd a login for the admin page that would determine what the user can do to the database."
